# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 08:55"

# --- Update Israel row (row 23) stats ---
$ws.Range("B23").Value = 3865
$ws.Range("C23").Value = 246
$ws.Range("E23").Value = 3764
$ws.Range("F23").Value = 66

# --- Kazajistan moves up the ranking (inserted between Bosnia y Herzegovina and
#     Jordania), pushing Jordania, Republica de Macedonia, Kuwait and Moldavia
#     down one row each; San Marino (row 86) is unaffected. ---
$ws.Range("A81").Value = "Kazajistan"
$ws.Range("B81").Value = 251
$ws.Range("C81").Value = 23
$ws.Range("D81").Value = 18
$ws.Range("E81").Value = 232
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 1

$ws.Range("A82").Value = "Jordania"
$ws.Range("B82").Value = 246
$ws.Range("C82").Value = 0
$ws.Range("D82").Value = 18
$ws.Range("E82").Value = 227
$ws.Range("F82").Value = 3
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 1

$ws.Range("A83").Value = "Republica de Macedonia"
$ws.Range("B83").Value = 241
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 3
$ws.Range("E83").Value = 234
$ws.Range("F83").Value = 1
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 4

$ws.Range("A84").Value = "Kuwait"
$ws.Range("B84").Value = 235
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 64
$ws.Range("E84").Value = 171
$ws.Range("F84").Value = 11
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0

$ws.Range("A85").Value = "Moldavia"
$ws.Range("B85").Value = 231
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 2
$ws.Range("E85").Value = 227
$ws.Range("F85").Value = 33
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 2

# --- Kirguistan moves up the ranking (inserted between Georgia and
#     Montenegro), pushing Montenegro, Bolivia, Trinidad y Tobago, Consejo
#     Danes para los Refugiados, Mayotte and Ruanda down one row each;
#     Liechtenstein (row 122) is unaffected. ---
$ws.Range("A115").Value = "Kirguistan"
$ws.Range("B115").Value = 84
$ws.Range("C115").Value = 26
$ws.Range("D115").Value = 0
$ws.Range("E115").Value = 84
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 0

$ws.Range("A116").Value = "Montenegro"
$ws.Range("B116").Value = 84
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 0
$ws.Range("E116").Value = 83
$ws.Range("F116").Value = 1
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 1

$ws.Range("A117").Value = "Bolivia"
$ws.Range("B117").Value = 81
$ws.Range("C117").Value = 7
$ws.Range("D117").Value = 0
$ws.Range("E117").Value = 81
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 0

$ws.Range("A118").Value = "Trinidad yTobago"
$ws.Range("B118").Value = 76
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 1
$ws.Range("E118").Value = 72
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 3

$ws.Range("A119").Value = "Consejo Danes para los Refugiados"
$ws.Range("B119").Value = 65
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 2
$ws.Range("E119").Value = 57
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 6

$ws.Range("A120").Value = "Mayotte"
$ws.Range("B120").Value = 63
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 0
$ws.Range("E120").Value = 63
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 0

$ws.Range("A121").Value = "Ruanda"
$ws.Range("B121").Value = 60
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 0
$ws.Range("E121").Value = 60
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 0
